$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034664657185821
$ws.Range("D2").Value = 1.035499414816682
$ws.Range("E2").Value = 1.043060927551002
$ws.Range("F2").Value = 1.052178916453339
$ws.Range("I2").Value = 1.034632786520877
$ws.Range("J2").Value = 1.039782661555336
$ws.Range("K2").Value = 1.038295843983805
$ws.Range("L2").Value = 1.045835857010109
$ws.Range("M2").Value = 1.054928366262623
$ws.Range("N2").Value = 1.041259272199483

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035609885579049
$ws.Range("D3").Value = 1.036313140980414
$ws.Range("E3").Value = 1.043907164975743
$ws.Range("F3").Value = 1.053123353420785
$ws.Range("I3").Value = 1.034758132182085
$ws.Range("J3").Value = 1.040371137546062
$ws.Range("K3").Value = 1.038918966212543
$ws.Range("L3").Value = 1.046492960411753
$ws.Range("M3").Value = 1.055685259188329
$ws.Range("N3").Value = 1.041848583893615

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03622198076301
$ws.Range("D4").Value = 1.036840392142673
$ws.Range("E4").Value = 1.044455525656433
$ws.Range("F4").Value = 1.053735337288572
$ws.Range("I4").Value = 1.034837906087498
$ws.Range("J4").Value = 1.040751774919679
$ws.Range("K4").Value = 1.039322233431118
$ws.Range("L4").Value = 1.046918293595544
$ws.Range("M4").Value = 1.056175263290503
$ws.Range("N4").Value = 1.042229761815962

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03647941686015
$ws.Range("D5").Value = 1.037062218522803
$ws.Range("E5").Value = 1.044686243839101
$ws.Range("F5").Value = 1.053992822135048
$ws.Range("I5").Value = 1.034871123347509
$ws.Range("J5").Value = 1.040911758898989
$ws.Range("K5").Value = 1.039491781541884
$ws.Range("L5").Value = 1.047097136846697
$ws.Range("M5").Value = 1.056381318102497
$ws.Range("N5").Value = 1.042389972990873

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036522648008893
$ws.Range("D6").Value = 1.037099474075748
$ws.Range("E6").Value = 1.044724993382601
$ws.Range("F6").Value = 1.054036067060142
$ws.Range("I6").Value = 1.034876681909227
$ws.Range("J6").Value = 1.0409386187736
$ws.Range("K6").Value = 1.039520250238231
$ws.Range("L6").Value = 1.047127167328233
$ws.Range("M6").Value = 1.056415918922827
$ws.Range("N6").Value = 1.042416871009586

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036225420200783
$ws.Range("D7").Value = 1.03684335553143
$ws.Range("E7").Value = 1.044458607789687
$ws.Range("F7").Value = 1.053738777002818
$ws.Range("I7").Value = 1.034838351194985
$ws.Range("J7").Value = 1.040753912776084
$ws.Range("K7").Value = 1.039324498885396
$ws.Range("L7").Value = 1.046920683178868
$ws.Range("M7").Value = 1.056178016382266
$ws.Range("N7").Value = 1.042231902708369

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034984003763312
$ws.Range("D8").Value = 1.03577426830451
$ws.Range("E8").Value = 1.043346753342027
$ws.Range("F8").Value = 1.052497912337111
$ws.Range("I8").Value = 1.03467542338114
$ws.Range("J8").Value = 1.039981569368429
$ws.Range("K8").Value = 1.038506416840226
$ws.Range("L8").Value = 1.0460578978178
$ws.Range("M8").Value = 1.055184110532602
$ws.Range("N8").Value = 1.041458462484485

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032800109228964
$ws.Range("D9").Value = 1.033895945000247
$ws.Range("E9").Value = 1.041393624882096
$ws.Range("F9").Value = 1.050318077465679
$ws.Range("I9").Value = 1.03437813743848
$ws.Range("J9").Value = 1.03861952947699
$ws.Range("K9").Value = 1.037065401787282
$ws.Range("L9").Value = 1.044538707769335
$ws.Range("M9").Value = 1.053434644523944
$ws.Range("N9").Value = 1.040094488340172

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031346687922498
$ws.Range("D10").Value = 1.032647540406588
$ws.Range("E10").Value = 1.040095724672435
$ws.Range("F10").Value = 1.048869462814493
$ws.Range("I10").Value = 1.034173129920349
$ws.Range("J10").Value = 1.037710842917555
$ws.Range("K10").Value = 1.03610515984171
$ws.Range("L10").Value = 1.04352675324782
$ws.Range("M10").Value = 1.052269701795692
$ws.Range("N10").Value = 1.039184511341595

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.03071794790163
$ws.Range("D11").Value = 1.032107888611014
$ws.Range("E11").Value = 1.039534729661128
$ws.Range("F11").Value = 1.04824330757773
$ws.Range("I11").Value = 1.034082749137727
$ws.Range("J11").Value = 1.037317227698308
$ws.Range("K11").Value = 1.03568948109659
$ws.Range("L11").Value = 1.043088778658581
$ws.Range("M11").Value = 1.051765608775338
$ws.Range("N11").Value = 1.038790337143589

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030484497112122
$ws.Range("D12").Value = 1.031907576842383
$ws.Range("E12").Value = 1.039326503407806
$ws.Range("F12").Value = 1.048010892985379
$ws.Range("I12").Value = 1.034048936134822
$ws.Range("J12").Value = 1.037171000240951
$ws.Range("K12").Value = 1.035535097444337
$ws.Range("L12").Value = 1.042926127800748
$ws.Range("M12").Value = 1.051578417719866
$ws.Range("N12").Value = 1.03864390202647

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030534568954335
$ws.Range("D13").Value = 1.03195053808753
$ws.Range("E13").Value = 1.039371161733175
$ws.Range("F13").Value = 1.048060739109172
$ws.Range("I13").Value = 1.034056200061532
$ws.Range("J13").Value = 1.037202367472189
$ws.Range("K13").Value = 1.035568212432855
$ws.Range("L13").Value = 1.042961015477862
$ws.Range("M13").Value = 1.05161856848483
$ws.Range("N13").Value = 1.038675313802774

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030698648919547
$ws.Range("D14").Value = 1.032091327935534
$ws.Range("E14").Value = 1.039517514499779
$ws.Range("F14").Value = 1.048224092693144
$ws.Range("I14").Value = 1.034079959069105
$ws.Range("J14").Value = 1.037305140918626
$ws.Range("K14").Value = 1.035676719323584
$ws.Range("L14").Value = 1.04307533322053
$ws.Range("M14").Value = 1.051750134445859
$ws.Range("N14").Value = 1.038778233199293

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030799756121484
$ws.Range("D15").Value = 1.032178091659009
$ws.Range("E15").Value = 1.039607707488386
$ws.Range("F15").Value = 1.048324762444196
$ws.Range("I15").Value = 1.034094565784416
$ws.Range("J15").Value = 1.037368460246083
$ws.Range("K15").Value = 1.035743576409377
$ws.Range("L15").Value = 1.043145772481588
$ws.Range("M15").Value = 1.051831203442551
$ws.Range("N15").Value = 1.038841642447458

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031388428070769
$ws.Range("D16").Value = 1.032683374727462
$ws.Range("E16").Value = 1.04013297734198
$ws.Range("F16").Value = 1.0489110420654
$ws.Range("I16").Value = 1.034179094292694
$ws.Range("J16").Value = 1.037736962830604
$ws.Range("K16").Value = 1.036132749528814
$ws.Range("L16").Value = 1.043555824663101
$ws.Range("M16").Value = 1.052303163949246
$ws.Range("N16").Value = 1.039210668347917

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031757848009219
$ws.Range("D17").Value = 1.033000571556778
$ws.Range("E17").Value = 1.040462735109879
$ws.Range("F17").Value = 1.04927909647343
$ws.Range("I17").Value = 1.034231685791555
$ws.Range("J17").Value = 1.037968075854568
$ws.Range("K17").Value = 1.036376898374636
$ws.Range("L17").Value = 1.043813096130431
$ws.Range("M17").Value = 1.05259930282556
$ws.Range("N17").Value = 1.039442109578884

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031973382332256
$ws.Range("D18").Value = 1.033185675353889
$ws.Range("E18").Value = 1.040655174048882
$ws.Range("F18").Value = 1.049493882853114
$ws.Range("I18").Value = 1.034262206034191
$ws.Range("J18").Value = 1.03810286583073
$ws.Range("K18").Value = 1.036519317035016
$ws.Range("L18").Value = 1.043963178258314
$ws.Range("M18").Value = 1.052772067849442
$ws.Range("N18").Value = 1.039577090972273

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032046883762238
$ws.Range("D19").Value = 1.033248805892682
$ws.Range("E19").Value = 1.040720807133013
$ws.Range("F19").Value = 1.049567137474418
$ws.Range("I19").Value = 1.034272586259266
$ws.Range("J19").Value = 1.038148823266875
$ws.Range("K19").Value = 1.036567879912139
$ws.Range("L19").Value = 1.044014355750173
$ws.Range("M19").Value = 1.052830981686993
$ws.Range("N19").Value = 1.039623113673248

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031718206747787
$ws.Range("D20").Value = 1.032966530200546
$ws.Range("E20").Value = 1.040427345198465
$ws.Range("F20").Value = 1.049239596695746
$ws.Range("I20").Value = 1.034226059300742
$ws.Range("J20").Value = 1.037943281090492
$ws.Range("K20").Value = 1.036350702392476
$ws.Range("L20").Value = 1.043785491256666
$ws.Range("M20").Value = 1.052567526584429
$ws.Range("N20").Value = 1.039417279603398

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030650328932759
$ws.Range("D21").Value = 1.032049864988045
$ws.Range("E21").Value = 1.039474413043211
$ws.Range("F21").Value = 1.048175984503362
$ws.Range("I21").Value = 1.034072969297623
$ws.Range("J21").Value = 1.037274877278581
$ws.Range("K21").Value = 1.035644766244915
$ws.Range("L21").Value = 1.043041668591252
$ws.Range("M21").Value = 1.051711390116997
$ws.Range("N21").Value = 1.038747926581408

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029979440181281
$ws.Range("D22").Value = 1.031474325708462
$ws.Range("E22").Value = 1.03887614845913
$ws.Range("F22").Value = 1.047508218191863
$ws.Range("I22").Value = 1.03397531821844
$ws.Range("J22").Value = 1.03685450287924
$ws.Range("K22").Value = 1.035201020502138
$ws.Range("L22").Value = 1.042574185790362
$ws.Range("M22").Value = 1.051173401925616
$ws.Range("N22").Value = 1.038326955202203

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030335040593965
$ws.Range("D23").Value = 1.031779353252339
$ws.Range("E23").Value = 1.039193215800105
$ws.Range("F23").Value = 1.047862121437966
$ws.Range("I23").Value = 1.034027217183939
$ws.Range("J23").Value = 1.037077362573608
$ws.Range("K23").Value = 1.035436248272167
$ws.Range("L23").Value = 1.042821989078756
$ws.Range("M23").Value = 1.051458570913517
$ws.Range("N23").Value = 1.038550131382899

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031736118748573
$ws.Range("D24").Value = 1.032981911751587
$ws.Range("E24").Value = 1.040443336077099
$ws.Range("F24").Value = 1.049257444616191
$ws.Range("I24").Value = 1.034228602150232
$ws.Range("J24").Value = 1.037954484821497
$ws.Range("K24").Value = 1.036362539196045
$ws.Range("L24").Value = 1.043797964649084
$ws.Range("M24").Value = 1.052581884800539
$ws.Range("N24").Value = 1.039428499244986

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033364260687623
$ws.Range("D25").Value = 1.034380870955061
$ws.Range("E25").Value = 1.04189782379944
$ws.Range("F25").Value = 1.050880811006087
$ws.Range("I25").Value = 1.034456196836769
$ws.Range("J25").Value = 1.038971770219853
$ws.Range("K25").Value = 1.037437866846959
$ws.Range("L25").Value = 1.04493131186717
$ws.Range("M25").Value = 1.053886687805035
$ws.Range("N25").Value = 1.040447229305294
